$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two changed cell values in row 2
$ws.Range("A2").Value = "Y"
$ws.Range("B2").Value = "PLP_TC_01"

# Update the selected cell to match the saved selection in the diff
$ws.Range("H15").Select()
